# "Alguns conflitos na aula 6"
# The "Clinica" / "Protocolo" lookup table (E2:H7) is replaced by a "Pedido"
# (order) table: idPedido / Data / idCliente / idProduto, with the "Data"
# column holding real dates. A fourth data row is added, and two extra
# formatted-but-empty rows are appended below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (merged E2:H2): "Clinica" -> "Pedido" ---
$ws.Range("E2").Value = "Pedido"

# --- Column headers (row 3) ---
$ws.Range("E3").Value = "idPedido"
$ws.Range("F3").Value = "Data"
$ws.Range("G3").Value = "idCliente"
$ws.Range("H3").Value = "idProduto"

# --- Data rows 4-6: keep idPedido/idCliente/idProduto numbers, replace the
#     old "Protocolo" text codes (00001/00002/00003) with real dates ---
$ws.Range("F4").Value = 43470
$ws.Range("F5").Value = 43961
$ws.Range("F6").Value = 44073
$ws.Range("F4:F6").NumberFormat = "mm-dd-yy"

# --- New data row 7 for the Pedido table ---
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 44219
$ws.Range("F7").NumberFormat = "mm-dd-yy"
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 4

# --- Two new, otherwise-empty rows under the table, carrying the date format ---
$ws.Range("F8").NumberFormat = "mm-dd-yy"
$ws.Range("F9").NumberFormat = "mm-dd-yy"

# --- Selection / view: was scrolled with topLeftCell B1 and H19 selected;
#     now back at the top with C8 selected ---
$ws.Range("C8").Select()
